# Updates cryptos list values (price + 1h volume change) per upstream scrape.
# Cells are stored as literal text (inlineStr) even when the text looks like a
# number (e.g. "41.52"), so each write forces Text format first, then clears
# the format again afterwards so the cell style index is left untouched (0),
# matching the original workbook -- only the stored value actually changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.ClearFormats()
}

Set-TextValue "D2" "27.952.07"
Set-TextValue "E2" "  -1.01%  "
Set-TextValue "D3" "1.868.20"
Set-TextValue "E4" "  -0.19%  "
Set-TextValue "D5" "312.69"
Set-TextValue "E5" "  -0.52%  "
Set-TextValue "E6" "  -0.24%  "
Set-TextValue "E7" "  -1.05%  "
Set-TextValue "D8" "0.3821"
Set-TextValue "E8" "  -2.76%  "
Set-TextValue "D9" "0.08936"
Set-TextValue "E9" "  -7.33%  "
Set-TextValue "D10" "1.118"
Set-TextValue "E10" "  -1.99%  "
Set-TextValue "D11" "41.52"
Set-TextValue "E11" "  -1.14%  "
Set-TextValue "D12" "6.365"
Set-TextValue "E12" "  -0.99%  "
Set-TextValue "D13" "20.67"
Set-TextValue "E13" "  -1.33%  "
Set-TextValue "D14" "1.862.29"
Set-TextValue "E14" "  -3.13%  "
Set-TextValue "D15" "7.229"
Set-TextValue "E15" "  -1.36%  "
Set-TextValue "E16" "  -0.14%  "
Set-TextValue "D17" "0.00001100"
Set-TextValue "E17" "  -2.10%  "
Set-TextValue "D18" "91.04"
Set-TextValue "E18" "  -1.73%  "
Set-TextValue "D19" "0.06647"
Set-TextValue "E19" "  +0.14%  "
Set-TextValue "D20" "18.14"
Set-TextValue "E20" "  +0.82%  "
Set-TextValue "E21" "  -0.18%  "
Set-TextValue "D22" "6.111"
Set-TextValue "E22" "  -1.87%  "
Set-TextValue "D23" "27.975.01"
Set-TextValue "E23" "  -1.13%  "
Set-TextValue "E24" "  +1.57%  "
Set-TextValue "E25" "  -2.04%  "
Set-TextValue "D26" "2.088.34"
Set-TextValue "E26" "  -2.56%  "
Set-TextValue "D27" "2.488"
Set-TextValue "E27" "  -6.30%  "
Set-TextValue "D28" "157.74"
Set-TextValue "E28" "  -0.18%  "
Set-TextValue "D29" "20.69"
Set-TextValue "E29" "  -1.30%  "
Set-TextValue "D30" "126.19"
Set-TextValue "E30" "  -0.70%  "
Set-TextValue "D31" "0.1064"
Set-TextValue "E31" "  -0.16%  "
Set-TextValue "D32" "1.053"
Set-TextValue "E32" "  -3.79%  "
Set-TextValue "D33" "5.594"
Set-TextValue "E33" "  -0.83%  "
Set-TextValue "D34" "3.601"
Set-TextValue "E34" "  -0.69%  "
Set-TextValue "D35" "9.443"
Set-TextValue "E35" "  -2.24%  "
Set-TextValue "D36" "0.06577"
Set-TextValue "E36" "  -1.42%  "
Set-TextValue "E37" "  -1.04%  "
Set-TextValue "D38" "0.2181"
Set-TextValue "E38" "  -0.62%  "
Set-TextValue "D39" "1.281"
Set-TextValue "E39" "  -1.48%  "
Set-TextValue "D40" "1.203"
Set-TextValue "E40" "  -3.01%  "
Set-TextValue "B41" "TheSandbox"
Set-TextValue "C41" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D41" "0.6366"
Set-TextValue "E41" "  -0.39%  "
Set-TextValue "B42" "Aptos"
Set-TextValue "C42" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D42" "11.48"
Set-TextValue "E42" "  +0.01%  "
Set-TextValue "D43" "4.887"
Set-TextValue "E43" "  -2.30%  "
Set-TextValue "B45" "Decentraland"
Set-TextValue "C45" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D45" "0.5998"
Set-TextValue "E45" "  -0.56%  "
Set-TextValue "B46" "EnergySwap"
Set-TextValue "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "13.14"
Set-TextValue "E46" "  -2.41%  "
Set-TextValue "D47" "1.283"
Set-TextValue "E47" "  -0.42%  "
Set-TextValue "E48" "  -2.31%  "
Set-TextValue "D49" "1.232"
Set-TextValue "E49" "  +3.71%  "
Set-TextValue "D50" "1.990"
Set-TextValue "E50" "  -2.51%  "
Set-TextValue "D51" "120.63"
Set-TextValue "E51" "  -2.33%  "
